$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-06-20 Thursday" "2024-06-21 Friday"

Replace-Text "61×56=3416" "78×74=5772"
Replace-Text "77×59=4543" "74×55=4070"
Replace-Text "99×42=4158" "90×52=4680"
Replace-Text "20×95=1900" "30×26=780"
Replace-Text "90×63=5670" "63×60=3780"

Replace-Text "29×82=2378" "56×42=2352"
Replace-Text "21×93=1953" "78×36=2808"
Replace-Text "17×87=1479" "35×38=1330"
Replace-Text "44×98=4312" "39×45=1755"
Replace-Text "71×56=3976" "69×74=5106"

Replace-Text "12×55=660" "34×96=3264"
Replace-Text "33×58=1914" "84×94=7896"
Replace-Text "31×91=2821" "68×48=3264"
Replace-Text "27×92=2484" "47×38=1786"
Replace-Text "12×58=696" "17×56=952"

Replace-Text "82×31=2542" "68×84=5712"
Replace-Text "32×91=2912" "58×62=3596"
Replace-Text "49×69=3381" "40×34=1360"
Replace-Text "75×70=5250" "20×96=1920"
Replace-Text "30×19=570" "18×91=1638"

Replace-Text "35×76=2660" "14×23=322"
Replace-Text "99×90=8910" "26×36=936"
Replace-Text "46×71=3266" "76×41=3116"
Replace-Text "98×68=6664" "12×98=1176"
Replace-Text "89×85=7565" "86×83=7138"
